# Auto-generated Excel COM-interop script applying the commit
# "Otomatik güncelleme: 2025-06-14 02:28:32" changes to before.xlsx

$wb = $excel.ActiveWorkbook

# ---- Sheet "eskalasyon": escalation rows get a new 45822 (2025-06-14) entry
#      per category, shifting every category block down by one row ----
$ws = $wb.Worksheets.Item("eskalasyon")

# Extend formatting from row 18 down through the newly added rows 19-27
# (keeps the date style used in column A consistent)
$ws.Range("A18:G18").Copy()
$ws.Range("A19:G27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(2,1).Value = 45673
$ws.Cells.Item(2,2).Value = "Motorin UltraForce"
$ws.Cells.Item(2,3).Value = 49.23
$ws.Cells.Item(2,4).Value = 0.06374243733794294
$ws.Cells.Item(2,5).Value = 0.06374243733794294
$ws.Cells.Item(2,6).Value = "Servis Diyarbakır"
$ws.Cells.Item(2,7).Value = 0.05

$ws.Cells.Item(3,1).Value = 45756
$ws.Cells.Item(3,2).Value = "Motorin UltraForce"
$ws.Cells.Item(3,3).Value = 46.38
$ws.Cells.Item(3,4).Value = -0.0578915295551492
$ws.Cells.Item(3,5).Value = -0.0578915295551492
$ws.Cells.Item(3,6).Value = "Servis Diyarbakır"
$ws.Cells.Item(3,7).Value = 0.05

$ws.Cells.Item(4,1).Value = 45822
$ws.Cells.Item(4,2).Value = "Motorin UltraForce"
$ws.Cells.Item(4,3).Value = 49.96
$ws.Cells.Item(4,4).Value = 0.07718844329452357
$ws.Cells.Item(4,5).Value = 0.07718844329452357
$ws.Cells.Item(4,6).Value = "Servis Diyarbakır"
$ws.Cells.Item(4,7).Value = 0.05

$ws.Cells.Item(5,1).Value = 45673
$ws.Cells.Item(5,2).Value = "Motorin UltraForce"
$ws.Cells.Item(5,3).Value = 48.88
$ws.Cells.Item(5,4).Value = 0.05663640293990491
$ws.Cells.Item(5,5).Value = 0.05663640293990491
$ws.Cells.Item(5,6).Value = "Servis Kayseri"
$ws.Cells.Item(5,7).Value = 0.05

$ws.Cells.Item(6,1).Value = 45756
$ws.Cells.Item(6,2).Value = "Motorin UltraForce"
$ws.Cells.Item(6,3).Value = 46.03
$ws.Cells.Item(6,4).Value = -0.05830605564648117
$ws.Cells.Item(6,5).Value = -0.05830605564648117
$ws.Cells.Item(6,6).Value = "Servis Kayseri"
$ws.Cells.Item(6,7).Value = 0.05

$ws.Cells.Item(7,1).Value = 45822
$ws.Cells.Item(7,2).Value = "Motorin UltraForce"
$ws.Cells.Item(7,3).Value = 49.64
$ws.Cells.Item(7,4).Value = 0.07842711275255265
$ws.Cells.Item(7,5).Value = 0.07842711275255265
$ws.Cells.Item(7,6).Value = "Servis Kayseri"
$ws.Cells.Item(7,7).Value = 0.05

$ws.Cells.Item(8,1).Value = 45673
$ws.Cells.Item(8,2).Value = "Motorin UltraForce"
$ws.Cells.Item(8,3).Value = 48.57
$ws.Cells.Item(8,4).Value = 0.05678851174934718
$ws.Cells.Item(8,5).Value = 0.05678851174934718
$ws.Cells.Item(8,6).Value = "Servis Samsun"
$ws.Cells.Item(8,7).Value = 0.05

$ws.Cells.Item(9,1).Value = 45756
$ws.Cells.Item(9,2).Value = "Motorin UltraForce"
$ws.Cells.Item(9,3).Value = 45.66
$ws.Cells.Item(9,4).Value = -0.05991352686843743
$ws.Cells.Item(9,5).Value = -0.05991352686843743
$ws.Cells.Item(9,6).Value = "Servis Samsun"
$ws.Cells.Item(9,7).Value = 0.05

$ws.Cells.Item(10,1).Value = 45822
$ws.Cells.Item(10,2).Value = "Motorin UltraForce"
$ws.Cells.Item(10,3).Value = 49.24
$ws.Cells.Item(10,4).Value = 0.07840560665790641
$ws.Cells.Item(10,5).Value = 0.07840560665790641
$ws.Cells.Item(10,6).Value = "Servis Samsun"
$ws.Cells.Item(10,7).Value = 0.05

$ws.Cells.Item(11,1).Value = 45673
$ws.Cells.Item(11,2).Value = "Motorin UltraForce"
$ws.Cells.Item(11,3).Value = 47.24
$ws.Cells.Item(11,4).Value = 0.07023108291798819
$ws.Cells.Item(11,5).Value = 0.07023108291798819
$ws.Cells.Item(11,6).Value = "Spot Araç Anadolu Toplama"
$ws.Cells.Item(11,7).Value = 0.05

$ws.Cells.Item(12,1).Value = 45756
$ws.Cells.Item(12,2).Value = "Motorin UltraForce"
$ws.Cells.Item(12,3).Value = 44.26
$ws.Cells.Item(12,4).Value = -0.06308213378492811
$ws.Cells.Item(12,5).Value = -0.06308213378492811
$ws.Cells.Item(12,6).Value = "Spot Araç Anadolu Toplama"
$ws.Cells.Item(12,7).Value = 0.05

$ws.Cells.Item(13,1).Value = 45822
$ws.Cells.Item(13,2).Value = "Motorin UltraForce"
$ws.Cells.Item(13,3).Value = 47.76
$ws.Cells.Item(13,4).Value = 0.07907817442385912
$ws.Cells.Item(13,5).Value = 0.07907817442385912
$ws.Cells.Item(13,6).Value = "Spot Araç Anadolu Toplama"
$ws.Cells.Item(13,7).Value = 0.05

$ws.Cells.Item(14,1).Value = 45673
$ws.Cells.Item(14,2).Value = "Motorin UltraForce"
$ws.Cells.Item(14,3).Value = 47.24
$ws.Cells.Item(14,4).Value = 0.07023108291798819
$ws.Cells.Item(14,5).Value = 0.07023108291798819
$ws.Cells.Item(14,6).Value = "Spot Araç Avrupa&Anadolu"
$ws.Cells.Item(14,7).Value = 0.05

$ws.Cells.Item(15,1).Value = 45756
$ws.Cells.Item(15,2).Value = "Motorin UltraForce"
$ws.Cells.Item(15,3).Value = 44.26
$ws.Cells.Item(15,4).Value = -0.06308213378492811
$ws.Cells.Item(15,5).Value = -0.06308213378492811
$ws.Cells.Item(15,6).Value = "Spot Araç Avrupa&Anadolu"
$ws.Cells.Item(15,7).Value = 0.05

$ws.Cells.Item(16,1).Value = 45822
$ws.Cells.Item(16,2).Value = "Motorin UltraForce"
$ws.Cells.Item(16,3).Value = 47.76
$ws.Cells.Item(16,4).Value = 0.07907817442385912
$ws.Cells.Item(16,5).Value = 0.07907817442385912
$ws.Cells.Item(16,6).Value = "Spot Araç Avrupa&Anadolu"
$ws.Cells.Item(16,7).Value = 0.05

$ws.Cells.Item(17,1).Value = 45673
$ws.Cells.Item(17,2).Value = "Motorin UltraForce"
$ws.Cells.Item(17,3).Value = 47.24
$ws.Cells.Item(17,4).Value = 0.07023108291798819
$ws.Cells.Item(17,5).Value = 0.07023108291798819
$ws.Cells.Item(17,6).Value = "Spot Araç Teknosa"
$ws.Cells.Item(17,7).Value = 0.05

$ws.Cells.Item(18,1).Value = 45756
$ws.Cells.Item(18,2).Value = "Motorin UltraForce"
$ws.Cells.Item(18,3).Value = 44.26
$ws.Cells.Item(18,4).Value = -0.06308213378492811
$ws.Cells.Item(18,5).Value = -0.06308213378492811
$ws.Cells.Item(18,6).Value = "Spot Araç Teknosa"
$ws.Cells.Item(18,7).Value = 0.05

$ws.Cells.Item(19,1).Value = 45822
$ws.Cells.Item(19,2).Value = "Motorin UltraForce"
$ws.Cells.Item(19,3).Value = 47.76
$ws.Cells.Item(19,4).Value = 0.07907817442385912
$ws.Cells.Item(19,5).Value = 0.07907817442385912
$ws.Cells.Item(19,6).Value = "Spot Araç Teknosa"
$ws.Cells.Item(19,7).Value = 0.05

$ws.Cells.Item(20,1).Value = 45673
$ws.Cells.Item(20,2).Value = "Motorin UltraForce"
$ws.Cells.Item(20,3).Value = 47.24
$ws.Cells.Item(20,4).Value = 0.07023108291798819
$ws.Cells.Item(20,5).Value = 0.07023108291798819
$ws.Cells.Item(20,6).Value = "TL/Desi Avrupa Toplama"
$ws.Cells.Item(20,7).Value = 0.05

$ws.Cells.Item(21,1).Value = 45756
$ws.Cells.Item(21,2).Value = "Motorin UltraForce"
$ws.Cells.Item(21,3).Value = 44.26
$ws.Cells.Item(21,4).Value = -0.06308213378492811
$ws.Cells.Item(21,5).Value = -0.06308213378492811
$ws.Cells.Item(21,6).Value = "TL/Desi Avrupa Toplama"
$ws.Cells.Item(21,7).Value = 0.05

$ws.Cells.Item(22,1).Value = 45822
$ws.Cells.Item(22,2).Value = "Motorin UltraForce"
$ws.Cells.Item(22,3).Value = 47.76
$ws.Cells.Item(22,4).Value = 0.07907817442385912
$ws.Cells.Item(22,5).Value = 0.07907817442385912
$ws.Cells.Item(22,6).Value = "TL/Desi Avrupa Toplama"
$ws.Cells.Item(22,7).Value = 0.05

$ws.Cells.Item(23,1).Value = 45784
$ws.Cells.Item(23,2).Value = "Motorin UltraForce"
$ws.Cells.Item(23,3).Value = 44.2
$ws.Cells.Item(23,4).Value = -0.05089113162980452
$ws.Cells.Item(23,5).Value = -0.05089113162980452
$ws.Cells.Item(23,6).Value = "TL/Desi Avrupa İade Toplama"
$ws.Cells.Item(23,7).Value = 0.05

$ws.Cells.Item(24,1).Value = 45822
$ws.Cells.Item(24,2).Value = "Motorin UltraForce"
$ws.Cells.Item(24,3).Value = 47.76
$ws.Cells.Item(24,4).Value = 0.0805429864253393
$ws.Cells.Item(24,5).Value = 0.0805429864253393
$ws.Cells.Item(24,6).Value = "TL/Desi Avrupa İade Toplama"
$ws.Cells.Item(24,7).Value = 0.05

$ws.Cells.Item(25,1).Value = 45673
$ws.Cells.Item(25,2).Value = "Motorin UltraForce"
$ws.Cells.Item(25,3).Value = 47.24
$ws.Cells.Item(25,4).Value = 0.05587840858292359
$ws.Cells.Item(25,5).Value = 0.05587840858292359
$ws.Cells.Item(25,6).Value = "TL/Desi Avrupa&Anadolu Dağıtım"
$ws.Cells.Item(25,7).Value = 0.05

$ws.Cells.Item(26,1).Value = 45756
$ws.Cells.Item(26,2).Value = "Motorin UltraForce"
$ws.Cells.Item(26,3).Value = 44.26
$ws.Cells.Item(26,4).Value = -0.06308213378492811
$ws.Cells.Item(26,5).Value = -0.06308213378492811
$ws.Cells.Item(26,6).Value = "TL/Desi Avrupa&Anadolu Dağıtım"
$ws.Cells.Item(26,7).Value = 0.05

$ws.Cells.Item(27,1).Value = 45822
$ws.Cells.Item(27,2).Value = "Motorin UltraForce"
$ws.Cells.Item(27,3).Value = 47.76
$ws.Cells.Item(27,4).Value = 0.07907817442385912
$ws.Cells.Item(27,5).Value = 0.07907817442385912
$ws.Cells.Item(27,6).Value = "TL/Desi Avrupa&Anadolu Dağıtım"
$ws.Cells.Item(27,7).Value = 0.05

# ---- Sheet "durum": StartDate column moves from 2025-06-13 to 2025-06-14 ----
$wsDurum = $wb.Worksheets.Item("durum")
# Force text formatting first so Excel keeps the value as a literal string
# instead of auto-converting the date-looking text into a date serial number
$wsDurum.Range("B2:B7").NumberFormat = "@"
$wsDurum.Cells.Item(2,2).Value = "2025-06-14"
$wsDurum.Cells.Item(3,2).Value = "2025-06-14"
$wsDurum.Cells.Item(4,2).Value = "2025-06-14"
$wsDurum.Cells.Item(5,2).Value = "2025-06-14"
$wsDurum.Cells.Item(6,2).Value = "2025-06-14"
$wsDurum.Cells.Item(7,2).Value = "2025-06-14"

# ---- District price-history sheets: append the new 2025-06-14 (45822) row ----
$wsData = $wb.Worksheets.Item("934015")
$wsData.Range("A192:C192").Copy()
$wsData.Range("A193:C193").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsData.Cells.Item(193,1).Value = 45822
$wsData.Cells.Item(193,2).Value = "Motorin UltraForce"
$wsData.Cells.Item(193,3).Value = 47.76

$wsData = $wb.Worksheets.Item("065001")
$wsData.Range("A192:C192").Copy()
$wsData.Range("A193:C193").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsData.Cells.Item(193,1).Value = 45822
$wsData.Cells.Item(193,2).Value = "Motorin UltraForce"
$wsData.Cells.Item(193,3).Value = 50.01

$wsData = $wb.Worksheets.Item("035001")
$wsData.Range("A192:C192").Copy()
$wsData.Range("A193:C193").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsData.Cells.Item(193,1).Value = 45822
$wsData.Cells.Item(193,2).Value = "Motorin UltraForce"
$wsData.Cells.Item(193,3).Value = 48.96

$wsData = $wb.Worksheets.Item("055001")
$wsData.Range("A192:C192").Copy()
$wsData.Range("A193:C193").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsData.Cells.Item(193,1).Value = 45822
$wsData.Cells.Item(193,2).Value = "Motorin UltraForce"
$wsData.Cells.Item(193,3).Value = 49.24

$wsData = $wb.Worksheets.Item("021001")
$wsData.Range("A192:C192").Copy()
$wsData.Range("A193:C193").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsData.Cells.Item(193,1).Value = 45822
$wsData.Cells.Item(193,2).Value = "Motorin UltraForce"
$wsData.Cells.Item(193,3).Value = 49.96

$wsData = $wb.Worksheets.Item("038001")
$wsData.Range("A192:C192").Copy()
$wsData.Range("A193:C193").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsData.Cells.Item(193,1).Value = 45822
$wsData.Cells.Item(193,2).Value = "Motorin UltraForce"
$wsData.Cells.Item(193,3).Value = 49.64

